# Updated Excel to YAML utility for dynamic handling of columns
#
# This script highlights the cells that are already populated in the
# "cubes", "joins" and "dimensions" / "measures" example rows (green for
# cubes/joins/dimensions value rows, yellow for the measures row and the
# dimensions "description" column), and appends three new notes to the
# bottom of the "cubes" sheet calling out fields present in the sample
# workbook but missing from this template.

$wb = $excel.ActiveWorkbook

# BGR color ints (Excel COM Interior/Font.Color is 0x00BBGGRR)
$yellow = 65535     # RGB FFFF00
$green  = 5296274   # RGB 92D050

$ws1 = $wb.Worksheets.Item(1)   # cubes
$ws2 = $wb.Worksheets.Item(2)   # joins
$ws3 = $wb.Worksheets.Item(3)   # dimensions
$ws4 = $wb.Worksheets.Item(4)   # measures

# ---------------------------------------------------------------
# cubes (sheet1): mark the filled-in example row (row 2) green and
# add the "missing from sample" notes below the table.
# ---------------------------------------------------------------
$ws1.Range("B2").Interior.Color = $green
$ws1.Range("D2").Interior.Color = $green
$ws1.Range("E2").Interior.Color = $green
$ws1.Range("F2").Interior.Color = $green

$ws1.Range("A8").Value = "2 things which are in sample excel but not in here -> "
$ws1.Range("A9").Value = "1. dim_measure_flag "
$ws1.Range("A10").Value = "2. views_col_names"

$ws1.Columns.Item(1).ColumnWidth = 14.25

# ---------------------------------------------------------------
# joins (sheet2): mark the filled-in example row (row 2) green.
# ---------------------------------------------------------------
$ws2.Range("C2").Interior.Color = $green
$ws2.Range("D2").Interior.Color = $green
$ws2.Range("E2").Interior.Color = $green

# ---------------------------------------------------------------
# dimensions (sheet3): mark the filled-in example row (row 2) green,
# except the long description column which gets highlighted yellow.
# ---------------------------------------------------------------
$ws3.Range("A2").Interior.Color = $green
$ws3.Range("B2").Interior.Color = $green
$ws3.Range("C2").Interior.Color = $yellow
$ws3.Range("D2").Interior.Color = $green
$ws3.Range("E2").Interior.Color = $green
$ws3.Range("F2").Interior.Color = $green

$ws3.Columns.Item(1).ColumnWidth = 15.4

# ---------------------------------------------------------------
# measures (sheet4): mark the filled-in example row (row 2) yellow.
# ---------------------------------------------------------------
$ws4.Range("A2").Interior.Color = $yellow
$ws4.Range("B2").Interior.Color = $yellow
$ws4.Range("C2").Interior.Color = $yellow
$ws4.Range("D2").Interior.Color = $yellow
$ws4.Range("E2").Interior.Color = $yellow

# ---------------------------------------------------------------
# Restore per-sheet selections seen in the target workbook, leaving
# "cubes" as the active (tab-selected) sheet.
# ---------------------------------------------------------------
$ws2.Range("D18").Select()
$ws3.Range("A15").Select()
$ws4.Range("C24").Select()
$ws1.Range("A11").Select()
$ws1.Activate()
